$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row, ColumnA (Forecasted Consumption MW), ColumnB (Timestamp serial)
$data = @(
    @(2, 5120, 45919),
    @(3, 5080, 45919.01041666666),
    @(4, 5060, 45919.02083333334),
    @(5, 5030, 45919.03125),
    @(6, 4980, 45919.04166666666),
    @(7, 4950, 45919.05208333334),
    @(8, 4940, 45919.0625),
    @(9, 4930, 45919.07291666666),
    @(10, 4900, 45919.08333333334),
    @(11, 4900, 45919.09375),
    @(12, 4900, 45919.10416666666),
    @(13, 4900, 45919.11458333334),
    @(14, 4920, 45919.125),
    @(15, 4940, 45919.13541666666),
    @(16, 4960, 45919.14583333334),
    @(17, 4990, 45919.15625),
    @(18, 5040, 45919.16666666666),
    @(19, 5110, 45919.17708333334),
    @(20, 5180, 45919.1875),
    @(21, 5270, 45919.19791666666),
    @(22, 5380, 45919.20833333334),
    @(23, 5510, 45919.21875),
    @(24, 5630, 45919.22916666666),
    @(25, 5760, 45919.23958333334),
    @(26, 5950, 45919.25),
    @(27, 6070, 45919.26041666666),
    @(28, 6140, 45919.27083333334),
    @(29, 6170, 45919.28125),
    @(30, 6190, 45919.29166666666),
    @(31, 6170, 45919.30208333334),
    @(32, 6120, 45919.3125),
    @(33, 6050, 45919.32291666666),
    @(34, 5940, 45919.33333333334),
    @(35, 5820, 45919.34375),
    @(36, 5690, 45919.35416666666),
    @(37, 5540, 45919.36458333334),
    @(38, 5380, 45919.375),
    @(39, 5240, 45919.38541666666),
    @(40, 5110, 45919.39583333334),
    @(41, 5010, 45919.40625),
    @(42, 4910, 45919.41666666666),
    @(43, 4840, 45919.42708333334),
    @(44, 4790, 45919.4375),
    @(45, 4750, 45919.44791666666),
    @(46, 4710, 45919.45833333334),
    @(47, 4700, 45919.46875),
    @(48, 4700, 45919.47916666666),
    @(49, 4700, 45919.48958333334),
    @(50, 4700, 45919.5),
    @(51, 4700, 45919.51041666666),
    @(52, 4700, 45919.52083333334),
    @(53, 4710, 45919.53125),
    @(54, 4740, 45919.54166666666),
    @(55, 4780, 45919.55208333334),
    @(56, 4810, 45919.5625),
    @(57, 4860, 45919.57291666666),
    @(58, 4930, 45919.58333333334),
    @(59, 4990, 45919.59375),
    @(60, 5060, 45919.60416666666),
    @(61, 5140, 45919.61458333334),
    @(62, 5240, 45919.625),
    @(63, 5340, 45919.63541666666),
    @(64, 5440, 45919.64583333334),
    @(65, 5530, 45919.65625),
    @(66, 5630, 45919.66666666666),
    @(67, 5740, 45919.67708333334),
    @(68, 5860, 45919.6875),
    @(69, 5980, 45919.69791666666),
    @(70, 6120, 45919.70833333334),
    @(71, 6240, 45919.71875),
    @(72, 6350, 45919.72916666666),
    @(73, 6470, 45919.73958333334),
    @(74, 6550, 45919.75),
    @(75, 6660, 45919.76041666666),
    @(76, 6780, 45919.77083333334),
    @(77, 6890, 45919.78125),
    @(78, 7020, 45919.79166666666),
    @(79, 7100, 45919.80208333334),
    @(80, 7100, 45919.8125),
    @(81, 7060, 45919.82291666666),
    @(82, 6930, 45919.83333333334),
    @(83, 6780, 45919.84375),
    @(84, 6670, 45919.85416666666),
    @(85, 6540, 45919.86458333334),
    @(86, 6380, 45919.875),
    @(87, 6240, 45919.88541666666),
    @(88, 6070, 45919.89583333334),
    @(89, 5920, 45919.90625),
    @(90, 5810, 45919.91666666666),
    @(91, 5670, 45919.92708333334),
    @(92, 5550, 45919.9375),
    @(93, 5440, 45919.94791666666),
    @(94, 5270, 45919.95833333334),
    @(95, 5220, 45919.96875),
    @(96, 5180, 45919.97916666666),
    @(97, 5140, 45919.98958333334)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
}
